# Apply the updated crypto price / volume(1h) figures scraped on
# Tue Aug 15 19:39:53 UTC 2023 (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.190.41"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "1.825.33"
$ws.Range("E3").Value = "  -0.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.20%  "

# Row 5
$ws.Range("D5").Value = "'236.22"
$ws.Range("E5").Value = "  -1.63%  "

# Row 6
$ws.Range("D6").Value = "'0.6017"
$ws.Range("E6").Value = "  -4.18%  "

# Row 7
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("D8").Value = "'0.07130"
$ws.Range("E8").Value = "  -3.88%  "

# Row 9
$ws.Range("E9").Value = "  -2.99%  "

# Row 10
$ws.Range("D10").Value = "'24.06"
$ws.Range("E10").Value = "  -3.01%  "

# Row 11
$ws.Range("D11").Value = "'0.07646"
$ws.Range("E11").Value = "  -1.11%  "

# Row 12
$ws.Range("D12").Value = "1.879.92"
$ws.Range("E12").Value = "  +2.09%  "

# Row 13
$ws.Range("D13").Value = "'4.769"
$ws.Range("E13").Value = "  -4.22%  "

# Row 14
$ws.Range("D14").Value = "'0.6400"
$ws.Range("E14").Value = "  -5.49%  "

# Row 15
$ws.Range("D15").Value = "'0.000009731"
$ws.Range("E15").Value = "  -4.45%  "

# Row 16
$ws.Range("D16").Value = "'79.29"
$ws.Range("E16").Value = "  -3.29%  "

# Row 17
$ws.Range("D17").Value = "2.045.60"
$ws.Range("E17").Value = "  -2.27%  "

# Row 18
$ws.Range("D18").Value = "'5.974"
$ws.Range("E18").Value = "  -4.34%  "

# Row 19
$ws.Range("D19").Value = "29.183.57"
$ws.Range("E19").Value = "  -0.61%  "

# Row 20
$ws.Range("D20").Value = "'230.66"
$ws.Range("E20").Value = "  +0.75%  "

# Row 21
$ws.Range("E21").Value = "  +0.24%  "

# Row 22
$ws.Range("D22").Value = "'11.68"
$ws.Range("E22").Value = "  -5.03%  "

# Row 23
$ws.Range("D23").Value = "'7.019"
$ws.Range("E23").Value = "  -5.34%  "

# Row 24
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("D25").Value = "'155.49"
$ws.Range("E25").Value = "  -2.23%  "

# Row 26
$ws.Range("D26").Value = "'8.029"
$ws.Range("E26").Value = "  -5.18%  "

# Row 27
$ws.Range("D27").Value = "'0.1275"
$ws.Range("E27").Value = "  -5.59%  "

# Row 28
$ws.Range("D28").Value = "'16.65"
$ws.Range("E28").Value = "  -4.44%  "

# Row 29
$ws.Range("D29").Value = "'0.06792"
$ws.Range("E29").Value = "  +5.32%  "

# Row 30
$ws.Range("D30").Value = "'1.452"
$ws.Range("E30").Value = "  +0.40%  "

# Row 31
$ws.Range("D31").Value = "'1.459"
$ws.Range("E31").Value = "  -1.78%  "

# Row 32
$ws.Range("D32").Value = "'3.791"
$ws.Range("E32").Value = "  -6.81%  "

# Row 33
$ws.Range("D33").Value = "'3.773"
$ws.Range("E33").Value = "  -7.20%  "

# Row 34
$ws.Range("D34").Value = "'1.132"
$ws.Range("E34").Value = "  -0.52%  "

# Row 35
$ws.Range("D35").Value = "'1.718"
$ws.Range("E35").Value = "  -6.40%  "

# Row 36
$ws.Range("D36").Value = "'0.6578"
$ws.Range("E36").Value = "  -4.86%  "

# Row 37
$ws.Range("E37").Value = "  -0.96%  "

# Row 38
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.230.97"
$ws.Range("E38").Value = "  -0.88%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.758"
$ws.Range("E39").Value = "  -2.28%  "

# Row 40
$ws.Range("D40").Value = "'0.01755"
$ws.Range("E40").Value = "  -5.52%  "

# Row 41
$ws.Range("D41").Value = "'6.505"
$ws.Range("E41").Value = "  -3.45%  "

# Row 42
$ws.Range("D42").Value = "'0.9237"
$ws.Range("E42").Value = "  -0.76%  "

# Row 43
$ws.Range("E43").Value = "  +0.27%  "

# Row 44
$ws.Range("D44").Value = "1.963.89"
$ws.Range("E44").Value = "  -2.94%  "

# Row 45
$ws.Range("D45").Value = "'99.78"
$ws.Range("E45").Value = "  -0.93%  "

# Row 46
$ws.Range("D46").Value = "'63.07"
$ws.Range("E46").Value = "  -3.95%  "

# Row 47
$ws.Range("D47").Value = "'0.00000000117"
$ws.Range("E47").Value = "  -1.38%  "

# Row 48
$ws.Range("D48").Value = "'1.621"
$ws.Range("E48").Value = "  -5.44%  "

# Row 49
$ws.Range("D49").Value = "'6.563"
$ws.Range("E49").Value = "  -6.95%  "

# Row 50
$ws.Range("D50").Value = "'0.05582"
$ws.Range("E50").Value = "  -1.67%  "

# Row 51
$ws.Range("D51").Value = "'8.444"
$ws.Range("E51").Value = "  -6.30%  "
